# GARP.xlsx data refresh: the data rows (A2:E31) need to be re-sorted in
# descending order by column D ("6mth return %"), which is how the sheet
# is presented. Row 1 (headers) stays in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:E31")
$sortKey   = $ws.Range("D2:D31")

$xlDescending = 2

$dataRange.Sort($sortKey, $xlDescending)
